$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2,  -0.3741840007719169,  0.1617489680992558),
    @(3,  -0.3707670419387925, -0.0556500001594112),
    @(4,  -0.2993579144914407, -0.03839939671659583),
    @(5,  -0.4597506293178769,  0.2066116146750451),
    @(6,  -0.2513257604274803, -0.2180752060329035),
    @(7,  -0.1857160465210193, -0.597420264083446),
    @(8,  -0.3404843653286578, -0.2748291975235543),
    @(9,  -0.3472996693999933,  0.07850268961801558),
    @(10, -0.1369763124766324,  0.3802984315508031),
    @(11, -0.1399738320131171,  0.1659283730374982),
    @(12, -0.01916413571477676, -0.2875422526141986),
    @(13, -0.01240185226818955, -0.2670929223789246),
    @(14,  0.1978784155580193, -0.232563517031517),
    @(15, -0.08103037563712231, 0.0628243135508008),
    @(16, -0.04720349089237135, -0.08339027837918121),
    @(17,  0.02245814608105134,  0.2219257256474202)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
